$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - write the values first
$headers = @("Job_Id","Job_Title","Job_Description","Total_Years_Min_Exp","Total_Years_Max_Exp","LinkedIn_Poster","LinkedIn_Posted","Resume_received","Resume_downloaded")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Build the header format (bold, centered/top aligned, thin box border) on
# A1 once, then clone it onto the rest of the header row via copy/paste of
# formats only, so every header cell shares a single combined cell style
# instead of each property write minting its own style record.
$a1 = $ws.Cells.Item(1, 1)
$a1.Font.Bold = $true
$a1.HorizontalAlignment = -4108
$a1.VerticalAlignment = -4160
$a1.Borders.LineStyle = 1

$a1.Copy()
$ws.Range("B1:I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data row
$ws.Cells.Item(2, 1).Value = "JD_001"
$ws.Cells.Item(2, 2).Value = "Senior Engineer"
$ws.Cells.Item(2, 3).Value = "We are seeking a Software Engineer to build and maintain high-quality software solutions.`nWork with global teams to drive innovation and deliver scalable applications.`nJoin Akkodis and be part of a tech-driven, collaborative environment."
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(2, 5).Value = 4

# The multi-line description auto-expands row 2's height; restore the
# default (non-custom) row height since wrap text was never enabled.
$ws.Rows.Item(2).AutoFit()
